# Fix some bug in Road Catalog
# Updates the condition expressions in column A (rows 8-13) of the
# C_Pedestrian.conf sheet so that the OR'd type comparisons are wrapped
# in parentheses and use the "%=%" comparison token instead of "=".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_Pedestrian.conf")

$ws.Range("A8").Value  = "facility_type=04;(type%=%00 || type%=%05 || type%=%07)"
$ws.Range("A9").Value  = "facility_type=04;(type%=%02 || type%=%03 || type%=%06)"
$ws.Range("A10").Value = "facility_type=03;(type%=%00 || type%=%05 || type%=%07)"
$ws.Range("A11").Value = "facility_type=03;(type%=%02 || type%=%03 || type%=%06)"
$ws.Range("A12").Value = "facility_type=05;(type%=%00 || type%=%05 || type%=%07)"
$ws.Range("A13").Value = "facility_type=05;(type%=%02 || type%=%03 || type%=%06)"

# Update the current selection shown when the file was last saved
# (the author had last selected B26:B27 before saving).
$ws.Range("B26:B27").Select()

# Update workbook window position recorded with the file.
$excel.Windows.Item(1).Left = 11250
$excel.Windows.Item(1).Top = 5115
